# "add UITip & bug fix"
# Adds three new localisation keys (wrong_character_selected_tip,
# no_character_selected_tip, select_character) with their Chinese values,
# and updates the sheet's selection/scroll state to point at the new cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header/key names (columns S, T, U) ---
$ws.Range("S1").Value = "wrong_character_selected_tip"
$ws.Range("T1").Value = "no_character_selected_tip"
$ws.Range("U1").Value = "select_character"

# --- Row 2: new Simplified Chinese translations (columns S, T, U) ---
$ws.Range("S2").Value = "错误：未知角色"
$ws.Range("T2").Value = "请选择一个角色再继续冒险吧！"
$ws.Range("U2").Value = "选择角色"

# --- Update the active selection / scroll position to match the edit ---
$ws.Range("T1").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$win.Height = 16760
